# ---------------------------------------------------------------------------
# Applies the Arabic-translation edits described by the target diff to the
# "ar_PS CrisisText Video Scripts for Demo.docx" document.
#
# Strategy:
#  - Plain paragraphs whose whole run-text changes are updated by assigning
#    directly to Paragraph.Range.Text (this bypasses Word's "smart quotes"
#    autocorrect that Find/Replace triggers, and keeps the paragraph's
#    rPr/pPr formatting intact).
#  - Paragraphs that contain a DOCPROPERTY field (so the run boundaries must
#    be preserved around the field) are updated by writing just the
#    sub-range of characters before and/or after the field.
#  - Paragraphs whose single run is split by manual line breaks (<w:br/>)
#    are recreated by assigning a single string that uses the vertical-tab
#    character (`v) as the line-break marker; Word's Range.Text setter turns
#    each `v into its own <w:t> segment joined by <w:br/>, mirroring the
#    original <w:t>/<w:br/> structure.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

function Set-ParagraphText {
    param($Doc, [int]$Index, [string]$Text)
    $Doc.Paragraphs.Item($Index).Range.Text = $Text
}

function Set-ParagraphAroundField {
    # Rewrites the text immediately before and/or after the (single) field
    # contained in paragraph $Index, leaving the field itself untouched.
    # Pass $null for $BeforeText / $AfterText to leave that side alone.
    param($Doc, [int]$Index, $BeforeText, $AfterText)

    $p = $Doc.Paragraphs.Item($Index)
    $pStart = $p.Range.Start
    $pEnd = $p.Range.End

    $fieldCodeStart = -1
    $fieldResultEnd = -1
    $fields = $Doc.Content.Fields
    for ($i = 1; $i -le $fields.Count; $i++) {
        $f = $fields.Item($i)
        $cs = $f.Code.Start
        $re = $f.Result.End
        if ($cs -ge $pStart -and $re -le $pEnd) {
            $fieldCodeStart = $cs
            $fieldResultEnd = $re
        }
    }
    if ($fieldCodeStart -eq -1) {
        throw "Set-ParagraphAroundField: no field found in paragraph $Index"
    }

    # Write the "after" text first so it cannot shift the "before" offsets.
    if ($AfterText -ne $null) {
        $afterRange = $Doc.Range($fieldResultEnd, $pEnd)
        $afterRange.Text = $AfterText
    }
    if ($BeforeText -ne $null) {
        $beforeRange = $Doc.Range($pStart, $fieldCodeStart)
        $beforeRange.Text = $BeforeText
    }
}

# --- 1. "What to expect" -> "ماذا نتوقع" --------------------------------
Set-ParagraphText $d 13 "ماذا نتوقع"

# --- 2. "Get more help." sentence: translate the English tail only -------
Set-ParagraphText $d 25 "للحصول على مزيد من المعلومات أو الموارد المتاحة لك في أوقات الأزمات، اختر“Get more help.” يمكنك أيضًا الوصول إلى هذه المعلومات عن طريق كتابة HELP في أي وقت. "

# --- 3. "Finally, selecting "Watch a video about {field}" will replay..." -
$finallyBefore = "أخيرًا، اختيار ""Watch a video about "
$finallyAfter = "`" سيُعيد تشغيل هذا الفيديو. "
Set-ParagraphAroundField $d 26 $finallyBefore $finallyAfter

# --- 4. MENU / "What would you like to do?" / Review Tips ----------------
$menuText = "القائمة" + "`v" + "“What would you like to do؟” " + "`v" + "`v" + "مراجعة النصائح"
Set-ParagraphText $d 27 $menuText

# --- 5. "Change my Settings" ---------------------------------------------
Set-ParagraphText $d 28 "تغيير الإعدادات"

# --- 6. "Invite a Friend to {field}" --------------------------------------
Set-ParagraphAroundField $d 29 "أدعو صديق ل " $null

# --- 7. "Get more help" ---------------------------------------------------
Set-ParagraphText $d 30 "احصل على المزيد من المساعدة"

# --- 8. "Watch a video about {field}" -------------------------------------
Set-ParagraphAroundField $d 31 "شاهد فيديو عن " $null

# --- 9. "Exit Menu" --------------------------------------------------------
Set-ParagraphText $d 32 "اخرج من القائمة"

# --- 10. "In very hard times, ... stability." -----------------------------
Set-ParagraphText $d 37 "في الأوقات الصعبة للغاية، ولا سيما خلال الحروب، قد يكون من الصعب أن نجد لحظات للتواصل مع أطفالنا، لكن هذه اللحظات، حتى وإن كانت صغيرة، هي التي يمكن أن تمنح أطفالنا الاستقرار الذي يحتا جونه بشدة. "

# --- 11. "{field} offers ideas on how to spend time..." -------------------
Set-ParagraphAroundField $d 39 $null "يقدم أفكارًا حول كيفية قضاء الوقت مع طفلك بواسطة الأنشطة المرحة. يمكنك القيام بهذه الأنشطة في أي مكان دون الحاجة إلى مستلزمات. "

# --- 12. "After each tip..." / "You can choose..." / "Active..." / "Calm..." -
$afterTipText = "بعد كل نصيحة، ستُسأل ما إذا كنت ترغب في تجرِبة نشاط مرحي أو إنهاء الدرس لليوم. " + "`v"
$afterTipText = $afterTipText + "يمكنك اختيار نوع النشاط الذي ترغب في القيام به مع أطفالك:" + "`v" + "`v"
$afterTipText = $afterTipText + "نشيط - متعة حركية" + "`v"
$afterTipText = $afterTipText + "هادى ء- للاسترخاء معًا، أو"
Set-ParagraphText $d 41 $afterTipText

# --- 13. "Quick - for when you are short on time" -------------------------
Set-ParagraphText $d 42 "سريع - عندما يكون لديك وقت محدود"

# --- 14. "You don't have to wait for the end of a tip..." -----------------
Set-ParagraphText $d 44 "لا يجب عليك الانتظار حتى نهاية النصيحة لتلقي هذه الأنشطة المرحة. يمكنك أيضًا كتابة PLAY في أي وقت. "
